# Updated symbol list on Mon Jan  2 20:47:21 UTC 2023 with GitHub Actions
# Refreshes the Price / Volume(1h) figures (and a few reordered coin rows)
# on the "cryptos" worksheet. Price/volume cells are stored as text, so we
# force the Text number format before writing the values to avoid Excel's
# automatic numeric/percentage conversion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "246.99"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.63%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "29.66"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "9.22%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.167"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.09%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05732"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.66%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.575"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.97%"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8571"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "4.58%"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8678"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.96%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1365"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.31%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.03396"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "5.89%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07077"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02929"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.71%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09381"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.14%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001514"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.65%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.04130"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2.07%"
$ws.Range("B16").Value = "One"
$ws.Range("C16").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0006016"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.15%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006083"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.08%"
$ws.Range("B18").Value = "UpBots"
$ws.Range("C18").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.007489"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "5,071.09%"
$ws.Range("B19").Value = "LEO"
$ws.Range("C19").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.489"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.63%"
$ws.Range("B20").Value = "GateToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.098"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.94%"
$ws.Range("B21").Value = "BTSEToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.186"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-5.68%"
$ws.Range("B22").Value = "BitpandaEcosystemToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.3184"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.63%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.26%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.462"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-3.06%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.41%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.001226"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.52%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001209"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "22.16%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03757"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.64%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.005736"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-3.74%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1074"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.42%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002448"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "1.80%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008496"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-12.55%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005245"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "2.13%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.05%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.06466"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-35.95%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002258"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-9.89%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.05%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.05%"